$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 tracks a rotating window of "blog" ticker entries (ser: 121/122/123)
# plus a meetup entry. The oldest blog entry (ser: 121) is retired and a new
# one (ser: 124) takes its place, with the remaining entries shifting up by
# one slot:
#   C8: blog ser:123 -> blog ser:124
#   E8: blog ser:122 -> blog ser:123
#   I8: blog ser:121 -> blog ser:122
# D8 (the meetup entry) is untouched.

$ws.Range("C8").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 124"
$ws.Range("E8").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 123"
$ws.Range("I8").Value = "type: blog" + [char]10 + "width: 2" + [char]10 + "height: 1" + [char]10 + "ser: 122"
